# Updated cryptos list on Fri May 10 05:49:39 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns for Sheet1 rows 2-51 to match
# the latest coinranking.com snapshot. Row 28 (Dai) is unchanged this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.740.41"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "3.033.20"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "594.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "152.65"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.027.52"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.517"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.65"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +12.55%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +1.65%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.57"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "3.534.76"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "62.745.22"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "3.031.53"
$ws.Range("E19").Value = "  +0.95%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "452.59"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  +1.36%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.10"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("E26").Value = "  +2.63%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.51%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.42"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +7.61%  "
$ws.Range("E32").Value = "  +0.09%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "27.56"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +8.76%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("E42").Value = "  -1.52%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.299"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +10.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "41.95"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.21%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "394.24"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "2.726.95"
$ws.Range("E47").Value = "  +0.38%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "132.21"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +0.04%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.97%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "24.38"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.44%  "
